$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing GDP per Capita values (South Korea, 1820-2010) ---
# Values are stored as text (shared strings), matching source data format.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "815.774301077617"
$ws.Range("E2").ClearFormats()
$ws.Range("E52").NumberFormat = "@"
$ws.Range("E52").Value = "820.357302769065"
$ws.Range("E52").ClearFormats()
$ws.Range("E93").NumberFormat = "@"
$ws.Range("E93").Value = "1107.55874209976"
$ws.Range("E93").ClearFormats()
$ws.Range("E94").NumberFormat = "@"
$ws.Range("E94").Value = "1090.4380672581"
$ws.Range("E94").ClearFormats()
$ws.Range("E95").NumberFormat = "@"
$ws.Range("E95").Value = "1171.04388664939"
$ws.Range("E95").ClearFormats()
$ws.Range("E96").NumberFormat = "@"
$ws.Range("E96").Value = "1254.40185371742"
$ws.Range("E96").ClearFormats()
$ws.Range("E97").NumberFormat = "@"
$ws.Range("E97").Value = "1261.26137963331"
$ws.Range("E97").ClearFormats()
$ws.Range("E98").NumberFormat = "@"
$ws.Range("E98").Value = "1314.99437366312"
$ws.Range("E98").ClearFormats()
$ws.Range("E99").NumberFormat = "@"
$ws.Range("E99").Value = "1366.14655876651"
$ws.Range("E99").ClearFormats()
$ws.Range("E100").NumberFormat = "@"
$ws.Range("E100").Value = "1427.0303064877"
$ws.Range("E100").ClearFormats()
$ws.Range("E101").NumberFormat = "@"
$ws.Range("E101").Value = "1363.15211392193"
$ws.Range("E101").ClearFormats()
$ws.Range("E102").NumberFormat = "@"
$ws.Range("E102").Value = "1428.81642393374"
$ws.Range("E102").ClearFormats()
$ws.Range("E103").NumberFormat = "@"
$ws.Range("E103").Value = "1473.87914036513"
$ws.Range("E103").ClearFormats()
$ws.Range("E104").NumberFormat = "@"
$ws.Range("E104").Value = "1492.98111473602"
$ws.Range("E104").ClearFormats()
$ws.Range("E105").NumberFormat = "@"
$ws.Range("E105").Value = "1505.50227622661"
$ws.Range("E105").ClearFormats()
$ws.Range("E106").NumberFormat = "@"
$ws.Range("E106").Value = "1422.14929874806"
$ws.Range("E106").ClearFormats()
$ws.Range("E107").NumberFormat = "@"
$ws.Range("E107").Value = "1469.8881628887"
$ws.Range("E107").ClearFormats()
$ws.Range("E108").NumberFormat = "@"
$ws.Range("E108").Value = "1529.54379534206"
$ws.Range("E108").ClearFormats()
$ws.Range("E109").NumberFormat = "@"
$ws.Range("E109").Value = "1611.71969844718"
$ws.Range("E109").ClearFormats()
$ws.Range("E110").NumberFormat = "@"
$ws.Range("E110").Value = "1492.93697243396"
$ws.Range("E110").ClearFormats()
$ws.Range("E111").NumberFormat = "@"
$ws.Range("E111").Value = "1495.86427487404"
$ws.Range("E111").ClearFormats()
$ws.Range("E112").NumberFormat = "@"
$ws.Range("E112").Value = "1637.20979738464"
$ws.Range("E112").ClearFormats()
$ws.Range("E113").NumberFormat = "@"
$ws.Range("E113").Value = "1620.77742636312"
$ws.Range("E113").ClearFormats()
$ws.Range("E114").NumberFormat = "@"
$ws.Range("E114").Value = "1620.94456310234"
$ws.Range("E114").ClearFormats()
$ws.Range("E115").NumberFormat = "@"
$ws.Range("E115").Value = "1655.02532507323"
$ws.Range("E115").ClearFormats()
$ws.Range("E116").NumberFormat = "@"
$ws.Range("E116").Value = "1673.30316238728"
$ws.Range("E116").ClearFormats()
$ws.Range("E117").NumberFormat = "@"
$ws.Range("E117").Value = "1759.5970076149"
$ws.Range("E117").ClearFormats()
$ws.Range("E118").NumberFormat = "@"
$ws.Range("E118").Value = "1692.88664269553"
$ws.Range("E118").ClearFormats()
$ws.Range("E119").NumberFormat = "@"
$ws.Range("E119").Value = "2005.74912487929"
$ws.Range("E119").ClearFormats()
$ws.Range("E120").NumberFormat = "@"
$ws.Range("E120").Value = "1935.09637809883"
$ws.Range("E120").ClearFormats()
$ws.Range("E121").NumberFormat = "@"
$ws.Range("E121").Value = "1610.73119637171"
$ws.Range("E121").ClearFormats()
$ws.Range("E122").NumberFormat = "@"
$ws.Range("E122").Value = "1975.23613458624"
$ws.Range("E122").ClearFormats()
$ws.Range("E123").NumberFormat = "@"
$ws.Range("E123").Value = "2107.24748381484"
$ws.Range("E123").ClearFormats()
$ws.Range("E124").NumberFormat = "@"
$ws.Range("E124").Value = "1645.72199188762"
$ws.Range("E124").ClearFormats()
$ws.Range("E125").NumberFormat = "@"
$ws.Range("E125").Value = "1685.93874842345"
$ws.Range("E125").ClearFormats()
$ws.Range("E128").NumberFormat = "@"
$ws.Range("E128").Value = "915.765104808997"
$ws.Range("E128").ClearFormats()
$ws.Range("E129").NumberFormat = "@"
$ws.Range("E129").Value = "937.262686713904"
$ws.Range("E129").ClearFormats()
$ws.Range("E130").NumberFormat = "@"
$ws.Range("E130").Value = "1098.87154619638"
$ws.Range("E130").ClearFormats()
$ws.Range("E131").NumberFormat = "@"
$ws.Range("E131").Value = "1188.98073473397"
$ws.Range("E131").ClearFormats()
$ws.Range("E132").NumberFormat = "@"
$ws.Range("E132").Value = "998.088764411796"
$ws.Range("E132").ClearFormats()
$ws.Range("E133").NumberFormat = "@"
$ws.Range("E133").Value = "990.770029025074"
$ws.Range("E133").ClearFormats()
$ws.Range("E134").NumberFormat = "@"
$ws.Range("E134").Value = "1038.04102619505"
$ws.Range("E134").ClearFormats()
$ws.Range("E135").NumberFormat = "@"
$ws.Range("E135").Value = "1316.69107643201"
$ws.Range("E135").ClearFormats()
$ws.Range("E136").NumberFormat = "@"
$ws.Range("E136").Value = "1372.67979924528"
$ws.Range("E136").ClearFormats()
$ws.Range("E137").NumberFormat = "@"
$ws.Range("E137").Value = "1411.33360203502"
$ws.Range("E137").ClearFormats()
$ws.Range("E138").NumberFormat = "@"
$ws.Range("E138").Value = "1381.98276854527"
$ws.Range("E138").ClearFormats()
$ws.Range("E139").NumberFormat = "@"
$ws.Range("E139").Value = "1466.76773051341"
$ws.Range("E139").ClearFormats()
$ws.Range("E140").NumberFormat = "@"
$ws.Range("E140").Value = "1517.88386369625"
$ws.Range("E140").ClearFormats()
$ws.Range("E141").NumberFormat = "@"
$ws.Range("E141").Value = "1555.72719655056"
$ws.Range("E141").ClearFormats()
$ws.Range("E142").NumberFormat = "@"
$ws.Range("E142").Value = "1547.69177263856"
$ws.Range("E142").ClearFormats()
$ws.Range("E143").NumberFormat = "@"
$ws.Range("E143").Value = "1605.84677683157"
$ws.Range("E143").ClearFormats()
$ws.Range("E144").NumberFormat = "@"
$ws.Range("E144").Value = "1620.49016551383"
$ws.Range("E144").ClearFormats()
$ws.Range("E145").NumberFormat = "@"
$ws.Range("E145").Value = "1720.82412067271"
$ws.Range("E145").ClearFormats()
$ws.Range("E146").NumberFormat = "@"
$ws.Range("E146").Value = "1834.9778125414"
$ws.Range("E146").ClearFormats()
$ws.Range("E147").NumberFormat = "@"
$ws.Range("E147").Value = "1917.42831867999"
$ws.Range("E147").ClearFormats()
$ws.Range("E148").NumberFormat = "@"
$ws.Range("E148").Value = "2093.88798132793"
$ws.Range("E148").ClearFormats()
$ws.Range("E149").NumberFormat = "@"
$ws.Range("E149").Value = "2232.05350605592"
$ws.Range("E149").ClearFormats()
$ws.Range("E150").NumberFormat = "@"
$ws.Range("E150").Value = "2468.48358458018"
$ws.Range("E150").ClearFormats()
$ws.Range("E151").NumberFormat = "@"
$ws.Range("E151").Value = "2764.15336626157"
$ws.Range("E151").ClearFormats()
$ws.Range("E152").NumberFormat = "@"
$ws.Range("E152").Value = "2974.80727329566"
$ws.Range("E152").ClearFormats()
$ws.Range("E153").NumberFormat = "@"
$ws.Range("E153").Value = "3221.67449377507"
$ws.Range("E153").ClearFormats()
$ws.Range("E154").NumberFormat = "@"
$ws.Range("E154").Value = "3387.89053461182"
$ws.Range("E154").ClearFormats()
$ws.Range("E155").NumberFormat = "@"
$ws.Range("E155").Value = "3822.04579946319"
$ws.Range("E155").ClearFormats()
$ws.Range("E156").NumberFormat = "@"
$ws.Range("E156").Value = "4112.60132503237"
$ws.Range("E156").ClearFormats()
$ws.Range("E157").NumberFormat = "@"
$ws.Range("E157").Value = "4362.00674670652"
$ws.Range("E157").ClearFormats()
$ws.Range("E158").NumberFormat = "@"
$ws.Range("E158").Value = "4855.94072488161"
$ws.Range("E158").ClearFormats()
$ws.Range("E159").NumberFormat = "@"
$ws.Range("E159").Value = "5367.79488735547"
$ws.Range("E159").ClearFormats()
$ws.Range("E160").NumberFormat = "@"
$ws.Range("E160").Value = "5856.49651833948"
$ws.Range("E160").ClearFormats()
$ws.Range("E161").NumberFormat = "@"
$ws.Range("E161").Value = "6265.88620524998"
$ws.Range("E161").ClearFormats()
$ws.Range("E162").NumberFormat = "@"
$ws.Range("E162").Value = "6064.04009282181"
$ws.Range("E162").ClearFormats()
$ws.Range("E163").NumberFormat = "@"
$ws.Range("E163").Value = "6398.8511360217"
$ws.Range("E163").ClearFormats()
$ws.Range("E164").NumberFormat = "@"
$ws.Range("E164").Value = "6821.47493608378"
$ws.Range("E164").ClearFormats()
$ws.Range("E165").NumberFormat = "@"
$ws.Range("E165").Value = "7611.73413986103"
$ws.Range("E165").ClearFormats()
$ws.Range("E166").NumberFormat = "@"
$ws.Range("E166").Value = "8303.51898630704"
$ws.Range("E166").ClearFormats()
$ws.Range("E167").NumberFormat = "@"
$ws.Range("E167").Value = "8859.35533634623"
$ws.Range("E167").ClearFormats()
$ws.Range("E168").NumberFormat = "@"
$ws.Range("E168").Value = "9756.2054506927"
$ws.Range("E168").ClearFormats()
$ws.Range("E169").NumberFormat = "@"
$ws.Range("E169").Value = "10864.9741316052"
$ws.Range("E169").ClearFormats()
$ws.Range("E170").NumberFormat = "@"
$ws.Range("E170").Value = "12039.9457868355"
$ws.Range("E170").ClearFormats()
$ws.Range("E171").NumberFormat = "@"
$ws.Range("E171").Value = "12759.4896120752"
$ws.Range("E171").ClearFormats()
$ws.Range("E172").NumberFormat = "@"
$ws.Range("E172").Value = "13874"
$ws.Range("E172").ClearFormats()
$ws.Range("E173").NumberFormat = "@"
$ws.Range("E173").Value = "15056.1603265677"
$ws.Range("E173").ClearFormats()
$ws.Range("E174").NumberFormat = "@"
$ws.Range("E174").Value = "15712.0022092603"
$ws.Range("E174").ClearFormats()
$ws.Range("E175").NumberFormat = "@"
$ws.Range("E175").Value = "16503.4338578176"
$ws.Range("E175").ClearFormats()
$ws.Range("E176").NumberFormat = "@"
$ws.Range("E176").Value = "17720.3953083733"
$ws.Range("E176").ClearFormats()
$ws.Range("E177").NumberFormat = "@"
$ws.Range("E177").Value = "19089.2991255748"
$ws.Range("E177").ClearFormats()
$ws.Range("E178").NumberFormat = "@"
$ws.Range("E178").Value = "20204.5613654713"
$ws.Range("E178").ClearFormats()
$ws.Range("E179").NumberFormat = "@"
$ws.Range("E179").Value = "21055.9461564496"
$ws.Range("E179").ClearFormats()
$ws.Range("E180").NumberFormat = "@"
$ws.Range("E180").Value = "19624.6915876868"
$ws.Range("E180").ClearFormats()
$ws.Range("E181").NumberFormat = "@"
$ws.Range("E181").Value = "21540.6925752225"
$ws.Range("E181").ClearFormats()
$ws.Range("E182").NumberFormat = "@"
$ws.Range("E182").Value = "23108.0175367617"
$ws.Range("E182").ClearFormats()
$ws.Range("E183").NumberFormat = "@"
$ws.Range("E183").Value = "23804.4376851571"
$ws.Range("E183").ClearFormats()
$ws.Range("E184").NumberFormat = "@"
$ws.Range("E184").Value = "25250.637414954"
$ws.Range("E184").ClearFormats()
$ws.Range("E185").NumberFormat = "@"
$ws.Range("E185").Value = "25679.2465612007"
$ws.Range("E185").ClearFormats()
$ws.Range("E186").NumberFormat = "@"
$ws.Range("E186").Value = "26646.5032864598"
$ws.Range("E186").ClearFormats()
$ws.Range("E187").NumberFormat = "@"
$ws.Range("E187").Value = "27443.283376883"
$ws.Range("E187").ClearFormats()
$ws.Range("E188").NumberFormat = "@"
$ws.Range("E188").Value = "28516.4837755568"
$ws.Range("E188").ClearFormats()
$ws.Range("E189").NumberFormat = "@"
$ws.Range("E189").Value = "29716.5704202778"
$ws.Range("E189").ClearFormats()
$ws.Range("E190").NumberFormat = "@"
$ws.Range("E190").Value = "30117.7285767672"
$ws.Range("E190").ClearFormats()
$ws.Range("E191").NumberFormat = "@"
$ws.Range("E191").Value = "29967.5630702629"
$ws.Range("E191").ClearFormats()
$ws.Range("E192").NumberFormat = "@"
$ws.Range("E192").Value = "31537.7729271793"
$ws.Range("E192").ClearFormats()

# --- Append new rows for years 2011-2016 ---
$ws.Range("A193").Value = 410
$ws.Range("B193").Value = "South Korea"
$ws.Range("C193").Value = "GDP per Capita"
$ws.Range("D193").Value = 2011
$ws.Range("E193").NumberFormat = "@"
$ws.Range("E193").Value = "32225"
$ws.Range("E193").ClearFormats()
$ws.Range("A194").Value = 410
$ws.Range("B194").Value = "South Korea"
$ws.Range("C194").Value = "GDP per Capita"
$ws.Range("D194").Value = 2012
$ws.Range("E194").NumberFormat = "@"
$ws.Range("E194").Value = "32791"
$ws.Range("E194").ClearFormats()
$ws.Range("A195").Value = 410
$ws.Range("B195").Value = "South Korea"
$ws.Range("C195").Value = "GDP per Capita"
$ws.Range("D195").Value = 2013
$ws.Range("E195").NumberFormat = "@"
$ws.Range("E195").Value = "33588"
$ws.Range("E195").ClearFormats()
$ws.Range("A196").Value = 410
$ws.Range("B196").Value = "South Korea"
$ws.Range("C196").Value = "GDP per Capita"
$ws.Range("D196").Value = 2014
$ws.Range("E196").NumberFormat = "@"
$ws.Range("E196").Value = "34493"
$ws.Range("E196").ClearFormats()
$ws.Range("A197").Value = 410
$ws.Range("B197").Value = "South Korea"
$ws.Range("C197").Value = "GDP per Capita"
$ws.Range("D197").Value = 2015
$ws.Range("E197").NumberFormat = "@"
$ws.Range("E197").Value = "35269"
$ws.Range("E197").ClearFormats()
$ws.Range("A198").Value = 410
$ws.Range("B198").Value = "South Korea"
$ws.Range("C198").Value = "GDP per Capita"
$ws.Range("D198").Value = 2016
$ws.Range("E198").NumberFormat = "@"
$ws.Range("E198").Value = "36103"
$ws.Range("E198").ClearFormats()
